# Terminplan.xlsx: "Version und Datum hinzugefügt"
# Adds a "Datum:" / date pair and a "Version:" / version-text pair
# below the existing "Bemerkungen" block on the Terminplan sheet
# (cells G20:H21).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Terminplan")

# xlRight = -4152, xlCenter = -4108
$xlRight  = -4152
$xlCenter = -4108

# --- Row 20: "Datum:" label + date value ----------------------------------
$ws.Range("G20").Value = "Datum:"
$ws.Range("G20").HorizontalAlignment = $xlRight
$ws.Range("G20").VerticalAlignment = $xlCenter

$ws.Range("H20").NumberFormat = "mm-dd-yy"
$ws.Range("H20").Value2 = 43766

# --- Row 21: "Version:" label + version text -------------------------------
$ws.Range("G21").Value = "Version:"
$ws.Range("G21").HorizontalAlignment = $xlRight
$ws.Range("G21").VerticalAlignment = $xlCenter

$ws.Range("H21").NumberFormat = "@"
$ws.Range("H21").Value = "1.1"
$ws.Range("H21").HorizontalAlignment = $xlRight
$ws.Range("H21").VerticalAlignment = $xlCenter

# Leave the selection where the author's last edit was made
$ws.Range("G21").Select() | Out-Null
